$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2049.3333
$ws.Range("J17").Value = 2049.3333
$ws.Range("L17").Value = 6147.999899999999
$ws.Range("N17").Value = -6483.999899999999
$ws.Range("H51").Value = 9312.357
$ws.Range("I51").Value = 9559.299999999999
$ws.Range("K51").Value = 9559.299999999999
$ws.Range("M51").Value = -9075.299999999999
$ws.Range("H132").Value = 4747.6577
$ws.Range("I132").Value = 4747.6577
$ws.Range("K132").Value = 14242.9731
$ws.Range("M132").Value = -11712.9731
$ws.Range("H135").Value = 1806.8
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("H137").Value = 14256.125
$ws.Range("I137").Value = 2134
$ws.Range("J137").Value = 50622.5
$ws.Range("K137").Value = 6402
$ws.Range("L137").Value = 151867.5
$ws.Range("M137").Value = -3852
$ws.Range("N137").Value = -156967.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3821.5293
$ws.Range("I61").Value = 2127.889
$ws.Range("K61").Value = 2127.889
$ws.Range("M61").Value = -1915.889
$ws.Range("H63").Value = 3998.4285
$ws.Range("I63").Value = 3998.4285
$ws.Range("K63").Value = 3998.4285
$ws.Range("M63").Value = -3312.4285
$ws.Range("H66").Value = 3998.4285
$ws.Range("I66").Value = 3998.4285
$ws.Range("K66").Value = 19992.1425
$ws.Range("M66").Value = -16560.1425
$ws.Range("H88").Value = 3871.077
$ws.Range("I88").Value = 1973.6666
$ws.Range("J88").Value = 5497.4287
$ws.Range("K88").Value = 1973.6666
$ws.Range("L88").Value = 5497.4287
$ws.Range("M88").Value = -1567.6666
$ws.Range("N88").Value = -6309.4287
$ws.Range("H91").Value = 3871.077
$ws.Range("I91").Value = 1973.6666
$ws.Range("J91").Value = 5497.4287
$ws.Range("K91").Value = 1973.6666
$ws.Range("L91").Value = 5497.4287
$ws.Range("M91").Value = -569.6666
$ws.Range("N91").Value = -8305.4287
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H102").Value = 5160.6924
$ws.Range("I102").Value = 4919
$ws.Range("K102").Value = 4919
$ws.Range("M102").Value = -3297
$ws.Range("H122").Value = 4786.423
$ws.Range("I122").Value = 4523.591
$ws.Range("K122").Value = 13570.773
$ws.Range("M122").Value = -11120.773
$ws.Range("H132").Value = 2083.6956
$ws.Range("I132").Value = 1386.5853
$ws.Range("J132").Value = 7800
$ws.Range("K132").Value = 4159.7559
$ws.Range("L132").Value = 23400
$ws.Range("M132").Value = -1629.7559
$ws.Range("N132").Value = -28460
$ws.Range("H136").Value = 3821.5293
$ws.Range("I136").Value = 2127.889
$ws.Range("K136").Value = 6383.667
$ws.Range("M136").Value = -3833.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2437.5881
$ws.Range("I86").Value = 2368.4614
$ws.Range("J86").Value = 2662.25
$ws.Range("K86").Value = 2368.4614
$ws.Range("L86").Value = 2662.25
$ws.Range("M86").Value = -1245.4614
$ws.Range("N86").Value = -4908.25
$ws.Range("H89").Value = 2437.5881
$ws.Range("I89").Value = 2368.4614
$ws.Range("J89").Value = 2662.25
$ws.Range("K89").Value = 11842.307
$ws.Range("L89").Value = 13311.25
$ws.Range("M89").Value = -6226.307000000001
$ws.Range("N89").Value = -24543.25
$ws.Range("H94").Value = 111111340
$ws.Range("I94").Value = 125000250
$ws.Range("K94").Value = 125000250
$ws.Range("M94").Value = -124999799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6762.294
$ws.Range("I31").Value = 4954
$ws.Range("J31").Value = 10077.5
$ws.Range("K31").Value = 4954
$ws.Range("L31").Value = 10077.5
$ws.Range("M31").Value = -4659
$ws.Range("N31").Value = -10667.5
$ws.Range("H34").Value = 6762.294
$ws.Range("I34").Value = 4954
$ws.Range("J34").Value = 10077.5
$ws.Range("K34").Value = 4954
$ws.Range("L34").Value = 10077.5
$ws.Range("M34").Value = -4752
$ws.Range("N34").Value = -10481.5
$ws.Range("H60").Value = 13306.7
$ws.Range("J60").Value = 21880.6
$ws.Range("L60").Value = 21880.6
$ws.Range("N60").Value = -22902.6
$ws.Range("H132").Value = 3251.5
$ws.Range("I132").Value = 2253.4736
$ws.Range("K132").Value = 6760.4208
$ws.Range("M132").Value = -4230.4208
$ws.Range("H134").Value = 5984.591
$ws.Range("I134").Value = 5533.2
$ws.Range("K134").Value = 16599.6
$ws.Range("M134").Value = -14064.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1561.3334
$ws.Range("I2").Value = 43.625
$ws.Range("J2").Value = 3295.8572
$ws.Range("K2").Value = 261.75
$ws.Range("L2").Value = 19775.1432
$ws.Range("M2").Value = -148.75
$ws.Range("N2").Value = -20001.1432
$ws.Range("H117").Value = 1000
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("H139").Value = 1704.5714
$ws.Range("I139").Value = 1211.5294
$ws.Range("K139").Value = 3634.5882
$ws.Range("M139").Value = 1505.4118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 122446.7
$ws.Range("I70").Value = 137974.14
$ws.Range("J70").Value = 5991
$ws.Range("K70").Value = 137974.14
$ws.Range("L70").Value = 5991
$ws.Range("M70").Value = -137704.14
$ws.Range("N70").Value = -6531
$ws.Range("H73").Value = 122446.7
$ws.Range("I73").Value = 137974.14
$ws.Range("J73").Value = 5991
$ws.Range("K73").Value = 137974.14
$ws.Range("L73").Value = 5991
$ws.Range("M73").Value = -137038.14
$ws.Range("N73").Value = -7863
$ws.Range("H102").Value = 5101.12
$ws.Range("I102").Value = 654.1579
$ws.Range("K102").Value = 654.1579
$ws.Range("M102").Value = 967.8421
$ws.Range("H113").Value = 4793.4
$ws.Range("I113").Value = 4676.8237
$ws.Range("K113").Value = 4676.8237
$ws.Range("M113").Value = -2506.8237
$ws.Range("H132").Value = 3816
$ws.Range("I132").Value = 1143.0416
$ws.Range("J132").Value = 25199.666
$ws.Range("K132").Value = 3429.1248
$ws.Range("L132").Value = 75598.99800000001
$ws.Range("M132").Value = -899.1248000000001
$ws.Range("N132").Value = -80658.99800000001
$ws.Range("H133").Value = 113995.6
$ws.Range("J133").Value = 113995.6
$ws.Range("L133").Value = 113995.6
$ws.Range("N133").Value = -124115.6
$ws.Range("H134").Value = 94700
$ws.Range("J134").Value = 94700
$ws.Range("L134").Value = 284100
$ws.Range("N134").Value = -289170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3425.5293
$ws.Range("I7").Value = 3459.6428
$ws.Range("J7").Value = 3266.3333
$ws.Range("K7").Value = 3459.6428
$ws.Range("L7").Value = 3266.3333
$ws.Range("M7").Value = -3347.6428
$ws.Range("N7").Value = -3490.3333
$ws.Range("H40").Value = 5409.727
$ws.Range("I40").Value = 5322.5
$ws.Range("J40").Value = 5898.2
$ws.Range("K40").Value = 5322.5
$ws.Range("L40").Value = 5898.2
$ws.Range("M40").Value = -5186.5
$ws.Range("N40").Value = -6170.2
$ws.Range("H92").Value = 49999
$ws.Range("J92").Value = 49999
$ws.Range("L92").Value = 49999
$ws.Range("N92").Value = -54991
$ws.Range("H100").Value = 4060.1428
$ws.Range("I100").Value = 4191.6665
$ws.Range("J100").Value = 3731.3333
$ws.Range("K100").Value = 4191.6665
$ws.Range("L100").Value = 3731.3333
$ws.Range("M100").Value = -3650.6665
$ws.Range("N100").Value = -4813.3333
$ws.Range("H122").Value = 4282.846
$ws.Range("I122").Value = 4282.846
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12848.538
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10398.538
$ws.Range("H126").Value = 3425.5293
$ws.Range("I126").Value = 3459.6428
$ws.Range("J126").Value = 3266.3333
$ws.Range("K126").Value = 10378.9284
$ws.Range("L126").Value = 9798.999899999999
$ws.Range("M126").Value = -7908.928400000001
$ws.Range("N126").Value = -14738.9999
$ws.Range("H132").Value = 5044.129
$ws.Range("I132").Value = 4168.6
$ws.Range("J132").Value = 6636
$ws.Range("K132").Value = 12505.8
$ws.Range("L132").Value = 19908
$ws.Range("M132").Value = -9975.800000000001
$ws.Range("N132").Value = -24968

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 25000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 25000
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("M49").Value = 25000
$ws.Range("N49").Value = -25460
$ws.Range("H62").Value = 5126.241
$ws.Range("I62").Value = 4154.84
$ws.Range("J62").Value = 11197.5
$ws.Range("K62").Value = 4154.84
$ws.Range("L62").Value = 11197.5
$ws.Range("M62").Value = -3530.84
$ws.Range("N62").Value = -12445.5
$ws.Range("H65").Value = 5126.241
$ws.Range("I65").Value = 4154.84
$ws.Range("J65").Value = 11197.5
$ws.Range("K65").Value = 20774.2
$ws.Range("L65").Value = 55987.5
$ws.Range("M65").Value = -17654.2
$ws.Range("N65").Value = -62227.5
$ws.Range("H81").Value = 3403.625
$ws.Range("I81").Value = 2989.8572
$ws.Range("K81").Value = 5979.7144
$ws.Range("M81").Value = -4918.7144
$ws.Range("H84").Value = 3403.625
$ws.Range("I84").Value = 2989.8572
$ws.Range("K84").Value = 29898.572
$ws.Range("M84").Value = -24594.572
$ws.Range("H122").Value = 11367842
$ws.Range("J122").Value = 41669172
$ws.Range("L122").Value = 125007516
$ws.Range("N122").Value = -125012416
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H126").Value = 2292.5386
$ws.Range("I126").Value = 2257.5557
$ws.Range("J126").Value = 2371.25
$ws.Range("K126").Value = 6772.6671
$ws.Range("L126").Value = 7113.75
$ws.Range("M126").Value = -4302.6671
$ws.Range("N126").Value = -12053.75
$ws.Range("H132").Value = 2278
$ws.Range("I132").Value = 2436.6897
$ws.Range("K132").Value = 7310.0691
$ws.Range("M132").Value = -4780.0691
